$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row 152 with the stock code and stock name
$ws.Cells.Item(152, 1).Value = "150.05.0511.00001"
$ws.Cells.Item(152, 2).Value = "M5X10 YILDIZ HAVSABAS CIVATA"

# Set explicit column widths (character units) so the saved bestFit widths
# land on the values Excel computed for this sheet's content
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 34
$ws.Columns.Item(3).ColumnWidth = 7.1666666666666667
$ws.Columns.Item(4).ColumnWidth = 7.1666666666666667
$ws.Columns.Item(5).ColumnWidth = 4.8333333333333333
$ws.Columns.Item(6).ColumnWidth = 6.8333333333333333
$ws.Columns.Item(7).ColumnWidth = 7.1666666666666667
$ws.Columns.Item(8).ColumnWidth = 4.8333333333333333
$ws.Columns.Item(9).ColumnWidth = 6.8333333333333333
$ws.Columns.Item(10).ColumnWidth = 10.6666666666666667
$ws.Columns.Item(11).ColumnWidth = 9.6666666666666667

# Update the view so the new row is visible / selected
$ws.Application.ActiveWindow.ScrollRow = 145
$ws.Range("A152").Select() | Out-Null
